$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column E header/type/field rows, copying formatting from column D ---

# E2 header cell (style like D2) - new shared string -> index 16
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "ItemInfoWindow에 뜰 string"

# E6 value cell (style like D6) - new shared string -> index 17
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = "날이 잘 든 단검이다`n가벼워서 쉽게 휘두를 수 있을 것 같다"

# E5 value cell (style like D6, since D5 will be cleared) - new shared string -> index 18
$ws.Range("D6").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = "더미이다`nㅁㄴㅇㄹ"

# E4 type cell (style like D4) - new shared string -> index 19
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = "itemDescription"

# E3 type cell (style like D3) - reuses existing "string" shared string
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "string"

$excel.CutCopyMode = 0

# --- Remove old shortsword value from D5 (data moved conceptually to new description column) ---
$ws.Range("D5").Clear()

# --- Resize rows for wrapped two-line description text ---
$ws.Rows.Item(5).RowHeight = 33
$ws.Rows.Item(6).RowHeight = 33

# --- Widen column E to fit the new description text ---
$ws.Columns.Item(5).ColumnWidth = 34.8

# --- Update active selection ---
$null = $ws.Range("D5").Select()
